$d = $word.ActiveDocument

# Locate the "...responded to supervision to date?" question text (form field
# statusText / label). The answer to this FORMTEXT field lives in the very
# next paragraph, which currently holds a handful of placeholder space runs
# and needs to become the {{response_to_probation}} merge field.
$rng = $d.Content
$found = $rng.Find.Execute("responded to supervision to date?", $true, $false, `
    $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the 'responded to supervision to date?' question text"
}

# Figure out which paragraph holds the found text, then grab the paragraph
# right after it - that is the FORMTEXT answer paragraph.
$count = $d.Paragraphs.Count
$answerPara = $null
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $rng.Start -and $p.Range.End -ge $rng.End) {
        $answerPara = $d.Paragraphs.Item($i + 1)
        break
    }
}

if ($null -eq $answerPara) {
    throw "Could not locate the answer paragraph for 'responded to supervision to date?'"
}

# Replace the paragraph's content (everything up to, but not including, the
# trailing paragraph mark) with the merge-field placeholder text, preserving
# the run's existing character formatting (Arial / noProof).
$answerRange = $d.Range($answerPara.Range.Start, $answerPara.Range.End - 1)
$answerRange.Text = "{{response_to_probation}}"
